# Add a new "Metadata" worksheet (placed after Sheet1) that records the
# locale the test fixture should be loaded/evaluated with, and make it the
# active tab.
$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Add()
$meta.Name = "Metadata"

$meta.Range("A1").Value = "Locale"
$meta.Range("B1").Value = "en-GB"

# Worksheets.Add() drops the new sheet in front of everything else by
# default -- move it so it comes right after "Sheet1".
$meta.Move($null, $wb.Worksheets.Item("Sheet1"))

# Make "Metadata" the active/selected tab, with B2 as the current selection.
$meta = $wb.Worksheets.Item("Metadata")
$meta.Activate()
[void]$meta.Range("B2").Select()
